$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value2 = 255.5  # H4: 201.42857 -> 255.5
$ws.Cells.Item(4, 9).Value2 = 233  # I4: 161.4 -> 233
$ws.Cells.Item(4, 10).Value2 = 300.5  # J4: 301.5 -> 300.5
$ws.Cells.Item(4, 11).Value2 = 233  # K4: 161.4 -> 233
$ws.Cells.Item(4, 12).Value2 = 300.5  # L4: 301.5 -> 300.5
$ws.Cells.Item(4, 13).Value2 = -119  # M4: -47.40000000000001 -> -119
$ws.Cells.Item(4, 14).Value2 = -528.5  # N4: -529.5 -> -528.5
$ws.Cells.Item(8, 8).Value2 = 54.333332  # H8: 29.285715 -> 54.333332
$ws.Cells.Item(8, 9).Value2 = 35.2  # I8: 29.285715 -> 35.2
$ws.Cells.Item(8, 10).Value2 = 150  # J8: 0 -> 150
$ws.Cells.Item(8, 11).Value2 = 105.6  # K8: 87.857145 -> 105.6
$ws.Cells.Item(8, 12).Value2 = 450  # L8: 0 -> 450
$ws.Cells.Item(8, 13).Value2 = 33.39999999999999  # M8: 51.142855 -> 33.39999999999999
$ws.Cells.Item(8, 14).Value2 = -728  # N8: <MISSING> -> -728
$ws.Cells.Item(18, 8).Value2 = 0  # H18: 4001 -> 0
$ws.Cells.Item(18, 10).Value2 = 0  # J18: 4001 -> 0
$ws.Cells.Item(18, 12).Value2 = 0  # L18: 4001 -> 0
$ws.Cells.Item(18, 14).ClearContents()  # N18: -4569 -> (removed)
$ws.Cells.Item(51, 8).Value2 = 3136.9744  # H51: 3137.2778 -> 3136.9744
$ws.Cells.Item(51, 10).Value2 = 3410.923  # J51: 3494.2 -> 3410.923
$ws.Cells.Item(51, 12).Value2 = 3410.923  # L51: 3494.2 -> 3410.923
$ws.Cells.Item(51, 14).Value2 = -4378.923  # N51: -4462.2 -> -4378.923
$ws.Cells.Item(76, 8).Value2 = 3242.423  # H76: 3246.87 -> 3242.423
$ws.Cells.Item(76, 9).Value2 = 3119.6516  # I76: 3123.0308 -> 3119.6516
$ws.Cells.Item(76, 11).Value2 = 3119.6516  # K76: 3123.0308 -> 3119.6516
$ws.Cells.Item(76, 13).Value2 = -2804.6516  # M76: -2808.0308 -> -2804.6516
$ws.Cells.Item(79, 8).Value2 = 3242.423  # H79: 3246.87 -> 3242.423
$ws.Cells.Item(79, 9).Value2 = 3119.6516  # I79: 3123.0308 -> 3119.6516
$ws.Cells.Item(79, 11).Value2 = 3119.6516  # K79: 3123.0308 -> 3119.6516
$ws.Cells.Item(79, 13).Value2 = -2027.6516  # M79: -2031.0308 -> -2027.6516
$ws.Cells.Item(94, 8).Value2 = 1498.6666  # H94: 1499 -> 1498.6666
$ws.Cells.Item(94, 9).Value2 = 1498.6666  # I94: 1499 -> 1498.6666
$ws.Cells.Item(94, 11).Value2 = 1498.6666  # K94: 1499 -> 1498.6666
$ws.Cells.Item(94, 13).Value2 = -1047.6666  # M94: -1048 -> -1047.6666
$ws.Cells.Item(106, 8).Value2 = 650  # H106: 637.125 -> 650
$ws.Cells.Item(106, 9).Value2 = 650  # I106: 637.125 -> 650
$ws.Cells.Item(106, 11).Value2 = 650  # K106: 637.125 -> 650
$ws.Cells.Item(106, 13).Value2 = -19  # M106: -6.125 -> -19
$ws.Cells.Item(131, 8).Value2 = 3820.9333  # H131: 3425.5293 -> 3820.9333
$ws.Cells.Item(131, 9).Value2 = 1231.4  # I131: 1102.8334 -> 1231.4
$ws.Cells.Item(131, 11).Value2 = 3694.2  # K131: 3308.5002 -> 3694.2
$ws.Cells.Item(131, 13).Value2 = 1345.8  # M131: 1731.4998 -> 1345.8
$ws.Cells.Item(132, 8).Value2 = 7174.2104  # H132: 6284.386 -> 7174.2104
$ws.Cells.Item(132, 9).Value2 = 7151.892  # I132: 6244.4883 -> 7151.892
$ws.Cells.Item(132, 11).Value2 = 21455.676  # K132: 18733.4649 -> 21455.676
$ws.Cells.Item(132, 13).Value2 = -18925.676  # M132: -16203.4649 -> -18925.676
$ws.Cells.Item(137, 8).Value2 = 23262086  # H137: 25006666 -> 23262086
$ws.Cells.Item(137, 9).Value2 = 34484170  # I137: 38462990 -> 34484170
$ws.Cells.Item(137, 10).Value2 = 16346.071  # J137: 16346.214 -> 16346.071
$ws.Cells.Item(137, 11).Value2 = 103452510  # K137: 115388970 -> 103452510
$ws.Cells.Item(137, 12).Value2 = 49038.213  # L137: 49038.642 -> 49038.213
$ws.Cells.Item(137, 13).Value2 = -103449960  # M137: -115386420 -> -103449960
$ws.Cells.Item(137, 14).Value2 = -54138.213  # N137: -54138.642 -> -54138.213
$ws.Cells.Item(138, 8).Value2 = 1819.2106  # H138: 1801.566 -> 1819.2106
$ws.Cells.Item(138, 10).Value2 = 2068.7556  # J138: 2070.2927 -> 2068.7556
$ws.Cells.Item(138, 12).Value2 = 6206.266799999999  # L138: 6210.8781 -> 6206.266799999999
$ws.Cells.Item(138, 14).Value2 = -16486.2668  # N138: -16490.8781 -> -16486.2668
$ws.Cells.Item(141, 8).Value2 = 4400  # H141: 5080 -> 4400
$ws.Cells.Item(141, 9).Value2 = 5040  # I141: 6050 -> 5040
$ws.Cells.Item(141, 11).Value2 = 15120  # K141: 18150 -> 15120
$ws.Cells.Item(141, 13).Value2 = -9940  # M141: -12970 -> -9940

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value2 = 7983  # H28: 6649.6665 -> 7983
$ws.Cells.Item(28, 9).Value2 = 9500  # I28: 7500 -> 9500
$ws.Cells.Item(28, 11).Value2 = 9500  # K28: 7500 -> 9500
$ws.Cells.Item(28, 13).Value2 = -9308  # M28: -7308 -> -9308
$ws.Cells.Item(32, 8).Value2 = 121136.93  # H32: 122561.06 -> 121136.93
$ws.Cells.Item(32, 9).Value2 = 133662.7  # I32: 135420.3 -> 133662.7
$ws.Cells.Item(32, 11).Value2 = 133662.7  # K32: 135420.3 -> 133662.7
$ws.Cells.Item(32, 13).Value2 = -133375.7  # M32: -135133.3 -> -133375.7
$ws.Cells.Item(46, 8).Value2 = 9950  # H46: 0 -> 9950
$ws.Cells.Item(46, 9).Value2 = 9950  # I46: 0 -> 9950
$ws.Cells.Item(46, 11).Value2 = 9950  # K46: 0 -> 9950
$ws.Cells.Item(46, 13).Value2 = -9631  # M46: <MISSING> -> -9631
$ws.Cells.Item(61, 8).Value2 = 7148029.5  # H61: 8341747 -> 7148029.5
$ws.Cells.Item(61, 9).Value2 = 6041  # I61: 10329.667 -> 6041
$ws.Cells.Item(61, 10).Value2 = 25003000  # J61: 33336000 -> 25003000
$ws.Cells.Item(61, 11).Value2 = 6041  # K61: 10329.667 -> 6041
$ws.Cells.Item(61, 12).Value2 = 25003000  # L61: 33336000 -> 25003000
$ws.Cells.Item(61, 13).Value2 = -5829  # M61: -10117.667 -> -5829
$ws.Cells.Item(61, 14).Value2 = -25003424  # N61: -33336424 -> -25003424
$ws.Cells.Item(92, 8).Value2 = 39199.4  # H92: 0 -> 39199.4
$ws.Cells.Item(92, 10).Value2 = 39199.4  # J92: 0 -> 39199.4
$ws.Cells.Item(92, 12).Value2 = 39199.4  # L92: 0 -> 39199.4
$ws.Cells.Item(92, 14).Value2 = -44191.4  # N92: <MISSING> -> -44191.4
$ws.Cells.Item(94, 8).Value2 = 42144.375  # H94: 42981.5 -> 42144.375
$ws.Cells.Item(94, 10).Value2 = 42144.375  # J94: 42981.5 -> 42144.375
$ws.Cells.Item(94, 12).Value2 = 42144.375  # L94: 42981.5 -> 42144.375
$ws.Cells.Item(94, 14).Value2 = -43946.375  # N94: -44783.5 -> -43946.375
$ws.Cells.Item(97, 8).Value2 = 2548.9  # H97: 2276.5557 -> 2548.9
$ws.Cells.Item(97, 10).Value2 = 3447.5  # J97: 2930 -> 3447.5
$ws.Cells.Item(97, 12).Value2 = 3447.5  # L97: 2930 -> 3447.5
$ws.Cells.Item(97, 14).Value2 = -4439.5  # N97: -3922 -> -4439.5
$ws.Cells.Item(99, 8).Value2 = 7983  # H99: 6649.6665 -> 7983
$ws.Cells.Item(99, 9).Value2 = 9500  # I99: 7500 -> 9500
$ws.Cells.Item(99, 11).Value2 = 9500  # K99: 7500 -> 9500
$ws.Cells.Item(99, 13).Value2 = -6505  # M99: -4505 -> -6505
$ws.Cells.Item(132, 8).Value2 = 4169835  # H132: 5003180 -> 4169835
$ws.Cells.Item(132, 9).Value2 = 8335437  # I132: 12501600 -> 8335437
$ws.Cells.Item(132, 11).Value2 = 25006311  # K132: 37504800 -> 25006311
$ws.Cells.Item(132, 13).Value2 = -25003781  # M132: -37502270 -> -25003781
$ws.Cells.Item(136, 8).Value2 = 7148029.5  # H136: 8341747 -> 7148029.5
$ws.Cells.Item(136, 9).Value2 = 6041  # I136: 10329.667 -> 6041
$ws.Cells.Item(136, 10).Value2 = 25003000  # J136: 33336000 -> 25003000
$ws.Cells.Item(136, 11).Value2 = 18123  # K136: 30989.001 -> 18123
$ws.Cells.Item(136, 12).Value2 = 75009000  # L136: 100008000 -> 75009000
$ws.Cells.Item(136, 13).Value2 = -15573  # M136: -28439.001 -> -15573
$ws.Cells.Item(136, 14).Value2 = -75014100  # N136: -100013100 -> -75014100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value2 = 47667.668  # H106: 53500 -> 47667.668
$ws.Cells.Item(106, 10).Value2 = 47667.668  # J106: 53500 -> 47667.668
$ws.Cells.Item(106, 12).Value2 = 47667.668  # L106: 53500 -> 47667.668
$ws.Cells.Item(106, 14).Value2 = -50191.668  # N106: -56024 -> -50191.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value2 = 15561.4  # H17: 13051.167 -> 15561.4
$ws.Cells.Item(17, 9).Value2 = 14451.75  # I17: 11661.4 -> 14451.75
$ws.Cells.Item(17, 11).Value2 = 14451.75  # K17: 11661.4 -> 14451.75
$ws.Cells.Item(17, 13).Value2 = -14277.75  # M17: -11487.4 -> -14277.75
$ws.Cells.Item(25, 8).Value2 = 3603.25  # H25: 3302.6 -> 3603.25
$ws.Cells.Item(25, 9).Value2 = 3603.25  # I25: 3302.6 -> 3603.25
$ws.Cells.Item(25, 11).Value2 = 3603.25  # K25: 3302.6 -> 3603.25
$ws.Cells.Item(25, 13).Value2 = -3429.25  # M25: -3128.6 -> -3429.25
$ws.Cells.Item(31, 8).Value2 = 2528514.8  # H31: 2418666.2 -> 2528514.8
$ws.Cells.Item(31, 9).Value2 = 6176544  # I31: 5559089.5 -> 6176544
$ws.Cells.Item(31, 11).Value2 = 6176544  # K31: 5559089.5 -> 6176544
$ws.Cells.Item(31, 13).Value2 = -6176249  # M31: -5558794.5 -> -6176249
$ws.Cells.Item(34, 8).Value2 = 2528514.8  # H34: 2418666.2 -> 2528514.8
$ws.Cells.Item(34, 9).Value2 = 6176544  # I34: 5559089.5 -> 6176544
$ws.Cells.Item(34, 11).Value2 = 6176544  # K34: 5559089.5 -> 6176544
$ws.Cells.Item(34, 13).Value2 = -6176342  # M34: -5558887.5 -> -6176342
$ws.Cells.Item(41, 8).Value2 = 12023.6  # H41: 4029.5 -> 12023.6
$ws.Cells.Item(41, 9).Value2 = 3706  # I41: 4029.5 -> 3706
$ws.Cells.Item(41, 10).Value2 = 24500  # J41: 0 -> 24500
$ws.Cells.Item(41, 11).Value2 = 3706  # K41: 4029.5 -> 3706
$ws.Cells.Item(41, 12).Value2 = 24500  # L41: 0 -> 24500
$ws.Cells.Item(41, 13).Value2 = -3278  # M41: -3601.5 -> -3278
$ws.Cells.Item(41, 14).Value2 = -25356  # N41: <MISSING> -> -25356
$ws.Cells.Item(58, 8).Value2 = 5566353.5  # H58: 5218644 -> 5566353.5
$ws.Cells.Item(58, 9).Value2 = 4467.2  # I58: 4333.727 -> 4467.2
$ws.Cells.Item(58, 11).Value2 = 4467.2  # K58: 4333.727 -> 4467.2
$ws.Cells.Item(58, 13).Value2 = -4264.2  # M58: -4130.727 -> -4264.2
$ws.Cells.Item(86, 8).Value2 = 20896.875  # H86: 9317.714 -> 20896.875
$ws.Cells.Item(86, 9).Value2 = 30045  # I86: 10644.8 -> 30045
$ws.Cells.Item(86, 10).Value2 = 5650  # J86: 6000 -> 5650
$ws.Cells.Item(86, 11).Value2 = 30045  # K86: 10644.8 -> 30045
$ws.Cells.Item(86, 12).Value2 = 5650  # L86: 6000 -> 5650
$ws.Cells.Item(86, 13).Value2 = -28922  # M86: -9521.799999999999 -> -28922
$ws.Cells.Item(86, 14).Value2 = -7896  # N86: -8246 -> -7896
$ws.Cells.Item(89, 8).Value2 = 20896.875  # H89: 9317.714 -> 20896.875
$ws.Cells.Item(89, 9).Value2 = 30045  # I89: 10644.8 -> 30045
$ws.Cells.Item(89, 10).Value2 = 5650  # J89: 6000 -> 5650
$ws.Cells.Item(89, 11).Value2 = 150225  # K89: 53224 -> 150225
$ws.Cells.Item(89, 12).Value2 = 28250  # L89: 30000 -> 28250
$ws.Cells.Item(89, 13).Value2 = -144609  # M89: -47608 -> -144609
$ws.Cells.Item(89, 14).Value2 = -39482  # N89: -41232 -> -39482
$ws.Cells.Item(94, 8).Value2 = 6766.8  # H94: 1796.4445 -> 6766.8
$ws.Cells.Item(94, 9).Value2 = 51499.5  # I94: 2999 -> 51499.5
$ws.Cells.Item(94, 10).Value2 = 1796.5  # J94: 1725.7059 -> 1796.5
$ws.Cells.Item(94, 11).Value2 = 51499.5  # K94: 2999 -> 51499.5
$ws.Cells.Item(94, 12).Value2 = 1796.5  # L94: 1725.7059 -> 1796.5
$ws.Cells.Item(94, 13).Value2 = -51048.5  # M94: -2548 -> -51048.5
$ws.Cells.Item(94, 14).Value2 = -2698.5  # N94: -2627.7059 -> -2698.5
$ws.Cells.Item(105, 8).Value2 = 18741.5  # H105: 1678 -> 18741.5
$ws.Cells.Item(105, 9).Value2 = 21509.8  # I105: 872.5 -> 21509.8
$ws.Cells.Item(105, 11).Value2 = 21509.8  # K105: 872.5 -> 21509.8
$ws.Cells.Item(105, 13).Value2 = -19762.8  # M105: 874.5 -> -19762.8
$ws.Cells.Item(132, 8).Value2 = 3090.8462  # H132: 2640.611 -> 3090.8462
$ws.Cells.Item(132, 9).Value2 = 2666.8  # I132: 2267.8667 -> 2666.8
$ws.Cells.Item(132, 11).Value2 = 8000.400000000001  # K132: 6803.6001 -> 8000.400000000001
$ws.Cells.Item(132, 13).Value2 = -5470.400000000001  # M132: -4273.6001 -> -5470.400000000001
$ws.Cells.Item(134, 8).Value2 = 3739.5715  # H134: 3876.6843 -> 3739.5715
$ws.Cells.Item(134, 9).Value2 = 2968.4443  # I134: 3034.875 -> 2968.4443
$ws.Cells.Item(134, 11).Value2 = 8905.332900000001  # K134: 9104.625 -> 8905.332900000001
$ws.Cells.Item(134, 13).Value2 = -6370.332900000001  # M134: -6569.625 -> -6370.332900000001
$ws.Cells.Item(136, 8).Value2 = 5566353.5  # H136: 5218644 -> 5566353.5
$ws.Cells.Item(136, 9).Value2 = 4467.2  # I136: 4333.727 -> 4467.2
$ws.Cells.Item(136, 11).Value2 = 13401.6  # K136: 13001.181 -> 13401.6
$ws.Cells.Item(136, 13).Value2 = -10851.6  # M136: -10451.181 -> -10851.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value2 = 7388.1665  # H3: 6765.5713 -> 7388.1665
$ws.Cells.Item(3, 9).Value2 = 3582.25  # I3: 3471.8 -> 3582.25
$ws.Cells.Item(3, 11).Value2 = 10746.75  # K3: 10415.4 -> 10746.75
$ws.Cells.Item(3, 13).Value2 = -10634.75  # M3: -10303.4 -> -10634.75
$ws.Cells.Item(10, 8).Value2 = 170.09091  # H10: 230.75 -> 170.09091
$ws.Cells.Item(10, 9).Value2 = 170.09091  # I10: 230.75 -> 170.09091
$ws.Cells.Item(10, 11).Value2 = 510.27273  # K10: 692.25 -> 510.27273
$ws.Cells.Item(10, 13).Value2 = -371.27273  # M10: -553.25 -> -371.27273
$ws.Cells.Item(33, 8).Value2 = 807.3077  # H33: 534.7 -> 807.3077
$ws.Cells.Item(33, 9).Value2 = 32.555557  # I33: 38.555557 -> 32.555557
$ws.Cells.Item(33, 10).Value2 = 2550.5  # J33: 5000 -> 2550.5
$ws.Cells.Item(33, 11).Value2 = 195.333342  # K33: 231.333342 -> 195.333342
$ws.Cells.Item(33, 12).Value2 = 15303  # L33: 30000 -> 15303
$ws.Cells.Item(33, 13).Value2 = 87.66665799999998  # M33: 51.66665799999998 -> 87.66665799999998
$ws.Cells.Item(33, 14).Value2 = -15869  # N33: -30566 -> -15869
$ws.Cells.Item(50, 8).Value2 = 2150.5217  # H50: 2161.9092 -> 2150.5217
$ws.Cells.Item(50, 9).Value2 = 2018.2222  # I50: 2025.1765 -> 2018.2222
$ws.Cells.Item(50, 11).Value2 = 6054.6666  # K50: 6075.529500000001 -> 6054.6666
$ws.Cells.Item(50, 13).Value2 = -5573.6666  # M50: -5594.529500000001 -> -5573.6666
$ws.Cells.Item(53, 8).Value2 = 2150.5217  # H53: 2161.9092 -> 2150.5217
$ws.Cells.Item(53, 9).Value2 = 2018.2222  # I53: 2025.1765 -> 2018.2222
$ws.Cells.Item(53, 11).Value2 = 6054.6666  # K53: 6075.529500000001 -> 6054.6666
$ws.Cells.Item(53, 13).Value2 = -5573.6666  # M53: -5594.529500000001 -> -5573.6666
$ws.Cells.Item(75, 8).Value2 = 2504  # H75: 2506 -> 2504
$ws.Cells.Item(75, 10).Value2 = 2504  # J75: 2506 -> 2504
$ws.Cells.Item(75, 12).Value2 = 7512  # L75: 7518 -> 7512
$ws.Cells.Item(75, 14).Value2 = -9508  # N75: -9514 -> -9508
$ws.Cells.Item(78, 8).Value2 = 2504  # H78: 2506 -> 2504
$ws.Cells.Item(78, 10).Value2 = 2504  # J78: 2506 -> 2504
$ws.Cells.Item(78, 12).Value2 = 22536  # L78: 22554 -> 22536
$ws.Cells.Item(78, 14).Value2 = -32520  # N78: -32538 -> -32520
$ws.Cells.Item(98, 8).Value2 = 1757  # H98: 2183.5 -> 1757
$ws.Cells.Item(98, 9).Value2 = 1759.8  # I98: 2400.25 -> 1759.8
$ws.Cells.Item(98, 11).Value2 = 5279.4  # K98: 7200.75 -> 5279.4
$ws.Cells.Item(98, 13).Value2 = -3781.4  # M98: -5702.75 -> -3781.4
$ws.Cells.Item(129, 8).Value2 = 1940.3334  # H129: 2017.4 -> 1940.3334
$ws.Cells.Item(129, 9).Value2 = 1019  # I129: 1098.5555 -> 1019
$ws.Cells.Item(129, 10).Value2 = 2993.2856  # J129: 3395.6667 -> 2993.2856
$ws.Cells.Item(129, 11).Value2 = 3057  # K129: 3295.6665 -> 3057
$ws.Cells.Item(129, 12).Value2 = 8979.856800000001  # L129: 10187.0001 -> 8979.856800000001
$ws.Cells.Item(129, 13).Value2 = 1943  # M129: 1704.3335 -> 1943
$ws.Cells.Item(129, 14).Value2 = -18979.8568  # N129: -20187.0001 -> -18979.8568
$ws.Cells.Item(134, 8).Value2 = 3060.1667  # H134: 3060.0557 -> 3060.1667
$ws.Cells.Item(134, 9).Value2 = 1872.2  # I134: 1872.0667 -> 1872.2
$ws.Cells.Item(134, 11).Value2 = 5616.6  # K134: 5616.2001 -> 5616.6
$ws.Cells.Item(134, 13).Value2 = -546.6000000000004  # M134: -546.2001 -> -546.6000000000004
$ws.Cells.Item(137, 8).Value2 = 8319.75  # H137: 7673.1113 -> 8319.75
$ws.Cells.Item(137, 9).Value2 = 4139.75  # I137: 3811.8 -> 4139.75
$ws.Cells.Item(137, 11).Value2 = 12419.25  # K137: 11435.4 -> 12419.25
$ws.Cells.Item(137, 13).Value2 = -7319.25  # M137: -6335.400000000001 -> -7319.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value2 = 4011917.2  # H11: 4012250.5 -> 4011917.2
$ws.Cells.Item(11, 10).Value2 = 18002.334  # J11: 19335.666 -> 18002.334
$ws.Cells.Item(11, 12).Value2 = 18002.334  # L11: 19335.666 -> 18002.334
$ws.Cells.Item(11, 14).Value2 = -18280.334  # N11: -19613.666 -> -18280.334
$ws.Cells.Item(70, 8).Value2 = 22039.625  # H70: 25396 -> 22039.625
$ws.Cells.Item(70, 9).Value2 = 11968.294  # I70: 16111.889 -> 11968.294
$ws.Cells.Item(70, 10).Value2 = 46498.57  # J70: 53248.332 -> 46498.57
$ws.Cells.Item(70, 11).Value2 = 11968.294  # K70: 16111.889 -> 11968.294
$ws.Cells.Item(70, 12).Value2 = 46498.57  # L70: 53248.332 -> 46498.57
$ws.Cells.Item(70, 13).Value2 = -11698.294  # M70: -15841.889 -> -11698.294
$ws.Cells.Item(70, 14).Value2 = -47038.57  # N70: -53788.332 -> -47038.57
$ws.Cells.Item(73, 8).Value2 = 22039.625  # H73: 25396 -> 22039.625
$ws.Cells.Item(73, 9).Value2 = 11968.294  # I73: 16111.889 -> 11968.294
$ws.Cells.Item(73, 10).Value2 = 46498.57  # J73: 53248.332 -> 46498.57
$ws.Cells.Item(73, 11).Value2 = 11968.294  # K73: 16111.889 -> 11968.294
$ws.Cells.Item(73, 12).Value2 = 46498.57  # L73: 53248.332 -> 46498.57
$ws.Cells.Item(73, 13).Value2 = -11032.294  # M73: -15175.889 -> -11032.294
$ws.Cells.Item(73, 14).Value2 = -48370.57  # N73: -55120.332 -> -48370.57
$ws.Cells.Item(80, 8).Value2 = 1866.25  # H80: 1899.5454 -> 1866.25
$ws.Cells.Item(80, 9).Value2 = 1699.25  # I80: 1765.6666 -> 1699.25
$ws.Cells.Item(80, 11).Value2 = 1699.25  # K80: 1765.6666 -> 1699.25
$ws.Cells.Item(80, 13).Value2 = -701.25  # M80: -767.6666 -> -701.25
$ws.Cells.Item(83, 8).Value2 = 1866.25  # H83: 1899.5454 -> 1866.25
$ws.Cells.Item(83, 9).Value2 = 1699.25  # I83: 1765.6666 -> 1699.25
$ws.Cells.Item(83, 11).Value2 = 8496.25  # K83: 8828.333000000001 -> 8496.25
$ws.Cells.Item(83, 13).Value2 = -3504.25  # M83: -3836.333000000001 -> -3504.25
$ws.Cells.Item(102, 8).Value2 = 1841.1351  # H102: 1902.0857 -> 1841.1351
$ws.Cells.Item(102, 9).Value2 = 1969.4546  # I102: 2006.8125 -> 1969.4546
$ws.Cells.Item(102, 10).Value2 = 782.5  # J102: 785 -> 782.5
$ws.Cells.Item(102, 11).Value2 = 1969.4546  # K102: 2006.8125 -> 1969.4546
$ws.Cells.Item(102, 12).Value2 = 782.5  # L102: 785 -> 782.5
$ws.Cells.Item(102, 13).Value2 = -347.4546  # M102: -384.8125 -> -347.4546
$ws.Cells.Item(102, 14).Value2 = -4026.5  # N102: -4029 -> -4026.5
$ws.Cells.Item(113, 8).Value2 = 3686.7778  # H113: 3593.2632 -> 3686.7778
$ws.Cells.Item(113, 9).Value2 = 2468.8  # I113: 2418 -> 2468.8
$ws.Cells.Item(113, 11).Value2 = 2468.8  # K113: 2418 -> 2468.8
$ws.Cells.Item(113, 13).Value2 = -298.8000000000002  # M113: -248 -> -298.8000000000002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value2 = 4.25  # H3: 2403.4 -> 4.25
$ws.Cells.Item(3, 10).Value2 = 0  # J3: 12000 -> 0
$ws.Cells.Item(3, 12).Value2 = 0  # L3: 12000 -> 0
$ws.Cells.Item(3, 14).ClearContents()  # N3: -12224 -> (removed)
$ws.Cells.Item(15, 8).Value2 = 4.25  # H15: 2403.4 -> 4.25
$ws.Cells.Item(15, 10).Value2 = 0  # J15: 12000 -> 0
$ws.Cells.Item(15, 12).Value2 = 0  # L15: 12000 -> 0
$ws.Cells.Item(15, 14).ClearContents()  # N15: -12340 -> (removed)
$ws.Cells.Item(22, 8).Value2 = 2287.0952  # H22: 2196.7273 -> 2287.0952
$ws.Cells.Item(22, 9).Value2 = 1022.8333  # I22: 919.4286 -> 1022.8333
$ws.Cells.Item(22, 11).Value2 = 1022.8333  # K22: 919.4286 -> 1022.8333
$ws.Cells.Item(22, 13).Value2 = -727.8333  # M22: -624.4286 -> -727.8333
$ws.Cells.Item(27, 8).Value2 = 2287.0952  # H27: 2196.7273 -> 2287.0952
$ws.Cells.Item(27, 9).Value2 = 1022.8333  # I27: 919.4286 -> 1022.8333
$ws.Cells.Item(27, 11).Value2 = 1022.8333  # K27: 919.4286 -> 1022.8333
$ws.Cells.Item(27, 13).Value2 = -915.8333  # M27: -812.4286 -> -915.8333
$ws.Cells.Item(40, 8).Value2 = 3894.7368  # H40: 4000.0557 -> 3894.7368
$ws.Cells.Item(40, 9).Value2 = 4036.5  # I40: 4193.231 -> 4036.5
$ws.Cells.Item(40, 11).Value2 = 4036.5  # K40: 4193.231 -> 4036.5
$ws.Cells.Item(40, 13).Value2 = -3900.5  # M40: -4057.231 -> -3900.5
$ws.Cells.Item(55, 8).Value2 = 1836.2222  # H55: 1661.9656 -> 1836.2222
$ws.Cells.Item(55, 9).Value2 = 2072.182  # I55: 1869.8182 -> 2072.182
$ws.Cells.Item(55, 10).Value2 = 1674  # J55: 1534.9445 -> 1674
$ws.Cells.Item(55, 11).Value2 = 2072.182  # K55: 1869.8182 -> 2072.182
$ws.Cells.Item(55, 12).Value2 = 1674  # L55: 1534.9445 -> 1674
$ws.Cells.Item(55, 13).Value2 = -1899.182  # M55: -1696.8182 -> -1899.182
$ws.Cells.Item(55, 14).Value2 = -2020  # N55: -1880.9445 -> -2020
$ws.Cells.Item(61, 8).Value2 = 23199.75  # H61: 12080.0625 -> 23199.75
$ws.Cells.Item(61, 9).Value2 = 24266.666  # I61: 10305.571 -> 24266.666
$ws.Cells.Item(61, 10).Value2 = 19999  # J61: 24501.5 -> 19999
$ws.Cells.Item(61, 11).Value2 = 24266.666  # K61: 10305.571 -> 24266.666
$ws.Cells.Item(61, 12).Value2 = 19999  # L61: 24501.5 -> 19999
$ws.Cells.Item(61, 13).Value2 = -24064.666  # M61: -10103.571 -> -24064.666
$ws.Cells.Item(61, 14).Value2 = -20403  # N61: -24905.5 -> -20403
$ws.Cells.Item(82, 8).Value2 = 1640.2  # H82: 1216.1177 -> 1640.2
$ws.Cells.Item(82, 9).Value2 = 3949.5  # I82: 1453 -> 3949.5
$ws.Cells.Item(82, 10).Value2 = 1062.875  # J82: 1050.3 -> 1062.875
$ws.Cells.Item(82, 11).Value2 = 3949.5  # K82: 1453 -> 3949.5
$ws.Cells.Item(82, 12).Value2 = 1062.875  # L82: 1050.3 -> 1062.875
$ws.Cells.Item(82, 13).Value2 = -3588.5  # M82: -1092 -> -3588.5
$ws.Cells.Item(82, 14).Value2 = -1784.875  # N82: -1772.3 -> -1784.875
$ws.Cells.Item(85, 8).Value2 = 1640.2  # H85: 1216.1177 -> 1640.2
$ws.Cells.Item(85, 9).Value2 = 3949.5  # I85: 1453 -> 3949.5
$ws.Cells.Item(85, 10).Value2 = 1062.875  # J85: 1050.3 -> 1062.875
$ws.Cells.Item(85, 11).Value2 = 3949.5  # K85: 1453 -> 3949.5
$ws.Cells.Item(85, 12).Value2 = 1062.875  # L85: 1050.3 -> 1062.875
$ws.Cells.Item(85, 13).Value2 = -2701.5  # M85: -205 -> -2701.5
$ws.Cells.Item(85, 14).Value2 = -3558.875  # N85: -3546.3 -> -3558.875
$ws.Cells.Item(93, 8).Value2 = 3703.3  # H93: 3307.9 -> 3703.3
$ws.Cells.Item(93, 9).Value2 = 2019.4286  # I93: 1898.25 -> 2019.4286
$ws.Cells.Item(93, 10).Value2 = 7632.3335  # J93: 8946.5 -> 7632.3335
$ws.Cells.Item(93, 11).Value2 = 2019.4286  # K93: 1898.25 -> 2019.4286
$ws.Cells.Item(93, 12).Value2 = 7632.3335  # L93: 8946.5 -> 7632.3335
$ws.Cells.Item(93, 13).Value2 = -771.4286  # M93: -650.25 -> -771.4286
$ws.Cells.Item(93, 14).Value2 = -10128.3335  # N93: -11442.5 -> -10128.3335
$ws.Cells.Item(113, 8).Value2 = 23199.75  # H113: 12080.0625 -> 23199.75
$ws.Cells.Item(113, 9).Value2 = 24266.666  # I113: 10305.571 -> 24266.666
$ws.Cells.Item(113, 10).Value2 = 19999  # J113: 24501.5 -> 19999
$ws.Cells.Item(113, 11).Value2 = 24266.666  # K113: 10305.571 -> 24266.666
$ws.Cells.Item(113, 12).Value2 = 19999  # L113: 24501.5 -> 19999
$ws.Cells.Item(113, 13).Value2 = -22096.666  # M113: -8135.571 -> -22096.666
$ws.Cells.Item(113, 14).Value2 = -24339  # N113: -28841.5 -> -24339
$ws.Cells.Item(123, 8).Value2 = 69999  # H123: 0 -> 69999
$ws.Cells.Item(123, 10).Value2 = 69999  # J123: 0 -> 69999
$ws.Cells.Item(123, 12).Value2 = 69999  # L123: 0 -> 69999
$ws.Cells.Item(123, 14).Value2 = -79799  # N123: <MISSING> -> -79799

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value2 = 412.66666  # H4: 286067.5 -> 412.66666
$ws.Cells.Item(4, 9).Value2 = 412.66666  # I4: 286067.5 -> 412.66666
$ws.Cells.Item(4, 11).Value2 = 412.66666  # K4: 286067.5 -> 412.66666
$ws.Cells.Item(4, 13).Value2 = -299.66666  # M4: -285954.5 -> -299.66666
$ws.Cells.Item(81, 8).Value2 = 1823.6111  # H81: 2107.818 -> 1823.6111
$ws.Cells.Item(81, 9).Value2 = 1548.5294  # I81: 1493.2632 -> 1548.5294
$ws.Cells.Item(81, 10).Value2 = 6500  # J81: 6000 -> 6500
$ws.Cells.Item(81, 11).Value2 = 3097.0588  # K81: 2986.5264 -> 3097.0588
$ws.Cells.Item(81, 12).Value2 = 13000  # L81: 12000 -> 13000
$ws.Cells.Item(81, 13).Value2 = -2036.0588  # M81: -1925.5264 -> -2036.0588
$ws.Cells.Item(81, 14).Value2 = -15122  # N81: -14122 -> -15122
$ws.Cells.Item(84, 8).Value2 = 1823.6111  # H84: 2107.818 -> 1823.6111
$ws.Cells.Item(84, 9).Value2 = 1548.5294  # I84: 1493.2632 -> 1548.5294
$ws.Cells.Item(84, 10).Value2 = 6500  # J84: 6000 -> 6500
$ws.Cells.Item(84, 11).Value2 = 15485.294  # K84: 14932.632 -> 15485.294
$ws.Cells.Item(84, 12).Value2 = 65000  # L84: 60000 -> 65000
$ws.Cells.Item(84, 13).Value2 = -10181.294  # M84: -9628.632000000001 -> -10181.294
$ws.Cells.Item(84, 14).Value2 = -75608  # N84: -70608 -> -75608
$ws.Cells.Item(113, 8).Value2 = 726.7  # H113: 559.19354 -> 726.7
$ws.Cells.Item(113, 9).Value2 = 922.4545000000001  # I113: 626.1 -> 922.4545000000001
$ws.Cells.Item(113, 10).Value2 = 487.44446  # J113: 437.54544 -> 487.44446
$ws.Cells.Item(113, 11).Value2 = 2767.3635  # K113: 1878.3 -> 2767.3635
$ws.Cells.Item(113, 12).Value2 = 1462.33338  # L113: 1312.63632 -> 1462.33338
$ws.Cells.Item(113, 13).Value2 = -597.3635000000004  # M113: 291.6999999999998 -> -597.3635000000004
$ws.Cells.Item(113, 14).Value2 = -5802.33338  # N113: -5652.63632 -> -5802.33338
$ws.Cells.Item(122, 8).Value2 = 57723.15  # H122: 71876.69 -> 57723.15
$ws.Cells.Item(122, 9).Value2 = 1553.1428  # I122: 1617.5834 -> 1553.1428
$ws.Cells.Item(122, 10).Value2 = 188786.5  # J122: 282654 -> 188786.5
$ws.Cells.Item(122, 11).Value2 = 4659.428400000001  # K122: 4852.7502 -> 4659.428400000001
$ws.Cells.Item(122, 12).Value2 = 566359.5  # L122: 847962 -> 566359.5
$ws.Cells.Item(122, 13).Value2 = -2209.428400000001  # M122: -2402.7502 -> -2209.428400000001
$ws.Cells.Item(122, 14).Value2 = -571259.5  # N122: -852862 -> -571259.5
$ws.Cells.Item(126, 8).Value2 = 2309.111  # H126: 2425.7778 -> 2309.111
$ws.Cells.Item(126, 9).Value2 = 2383.6428  # I126: 2497.7693 -> 2383.6428
$ws.Cells.Item(126, 10).Value2 = 2048.25  # J126: 2238.6 -> 2048.25
$ws.Cells.Item(126, 11).Value2 = 7150.928400000001  # K126: 7493.3079 -> 7150.928400000001
$ws.Cells.Item(126, 12).Value2 = 6144.75  # L126: 6715.799999999999 -> 6144.75
$ws.Cells.Item(126, 13).Value2 = -4680.928400000001  # M126: -5023.3079 -> -4680.928400000001
$ws.Cells.Item(126, 14).Value2 = -11084.75  # N126: -11655.8 -> -11084.75
$ws.Cells.Item(129, 8).Value2 = 64999.5  # H129: 67999.5 -> 64999.5
$ws.Cells.Item(129, 10).Value2 = 64999.5  # J129: 67999.5 -> 64999.5
$ws.Cells.Item(129, 12).Value2 = 64999.5  # L129: 67999.5 -> 64999.5
$ws.Cells.Item(129, 14).Value2 = -74999.5  # N129: -77999.5 -> -74999.5
$ws.Cells.Item(132, 8).Value2 = 5557176  # H132: 5748790 -> 5557176
$ws.Cells.Item(132, 10).Value2 = 1876.1666  # J132: 2176.6 -> 1876.1666
$ws.Cells.Item(132, 12).Value2 = 5628.4998  # L132: 6529.799999999999 -> 5628.4998
$ws.Cells.Item(132, 14).Value2 = -10688.4998  # N132: -11589.8 -> -10688.4998
$ws.Cells.Item(136, 8).Value2 = 14303282  # H136: 18276318 -> 14303282
$ws.Cells.Item(136, 9).Value2 = 6787794  # I136: 8597849 -> 6787794
$ws.Cells.Item(136, 10).Value2 = 50001850  # J136: 66668668 -> 50001850
$ws.Cells.Item(136, 11).Value2 = 20363382  # K136: 25793547 -> 20363382
$ws.Cells.Item(136, 12).Value2 = 150005550  # L136: 200006004 -> 150005550
$ws.Cells.Item(136, 13).Value2 = -20360832  # M136: -25790997 -> -20360832
$ws.Cells.Item(136, 14).Value2 = -150010650  # N136: -200011104 -> -150010650
